$d = $word.ActiveDocument

# This edit corresponds to regenerating the duplicatedBookmarks expected
# output with a newer engine run: the internal run-id used for the two
# "REF bookmark1" field runs changes, and the randomly generated id used
# for the "bookmark1" bookmarkStart/bookmarkEnd pair changes.

$oldRunRsid = "5BEF7D3A56B04B6B848A298613E2C005"
$newRunRsid = "57DEED68B53FB6D46541DB060630C897"

$oldBookmarkId = "113640858737380756001237403724904465710"
$newBookmarkId = "73557893092093159662871032680731948822"

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Paragraph 2: "Test link before bookmark : <REF bookmark1 field>"
$p2 = $d.Paragraphs.Item(2)
$p2Xml = '<w:p ' + $wNs + ' w:rsidP="009168BC" w:rsidR="00E02A2B" w:rsidRDefault="00E02A2B"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs></w:pPr>' +
         '<w:r><w:t xml:space="preserve">Test link before bookmark : </w:t></w:r>' +
         '<w:r w:rsidR="' + $newRunRsid + '"><w:fldChar w:fldCharType="begin"/></w:r>' +
         '<w:r w:rsidR="' + $newRunRsid + '"><w:instrText xml:space="preserve"> REF bookmark1 \h </w:instrText></w:r>' +
         '<w:r w:rsidR="' + $newRunRsid + '"><w:fldChar w:fldCharType="separate"/></w:r>' +
         '<w:r w:rsidR="' + $newRunRsid + '"><w:rPr><w:b w:val="true"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r>' +
         '<w:r w:rsidR="' + $newRunRsid + '"><w:fldChar w:fldCharType="end"/></w:r>' +
         '</w:p>'
$null = $p2.Range.InsertXML($p2Xml)

# Paragraph 3: "Test bookmark : " + bookmarkStart/bookmarkEnd pair around "bookmarked content"
$p3 = $d.Paragraphs.Item(3)
$p3Xml = '<w:p ' + $wNs + ' w:rsidP="00C31A62" w:rsidR="00C31A62" w:rsidRDefault="00C31A62"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs></w:pPr>' +
         '<w:r><w:t xml:space="preserve">Test bookmark : </w:t></w:r>' +
         '<w:bookmarkStart w:name="bookmark1" w:id="' + $newBookmarkId + '"/>' +
         '<w:r><w:t>bookmarked content</w:t></w:r>' +
         '<w:bookmarkEnd w:id="' + $newBookmarkId + '"/>' +
         '</w:p>'
$null = $p3.Range.InsertXML($p3Xml)

# Paragraph 5: "Test link after bookmark : <REF bookmark1 field>"
$p5 = $d.Paragraphs.Item(5)
$p5Xml = '<w:p ' + $wNs + ' w:rsidP="00E02A2B" w:rsidR="00E02A2B" w:rsidRDefault="00E02A2B"><w:pPr><w:tabs><w:tab w:pos="3119" w:val="left"/></w:tabs></w:pPr>' +
         '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
         '<w:r><w:t xml:space="preserve">Test link after bookmark : </w:t></w:r>' +
         '<w:r w:rsidR="' + $newRunRsid + '"><w:fldChar w:fldCharType="begin"/></w:r>' +
         '<w:r w:rsidR="' + $newRunRsid + '"><w:instrText xml:space="preserve"> REF bookmark1 \h </w:instrText></w:r>' +
         '<w:r w:rsidR="' + $newRunRsid + '"><w:fldChar w:fldCharType="separate"/></w:r>' +
         '<w:r w:rsidR="' + $newRunRsid + '"><w:rPr><w:b w:val="true"/><w:noProof/></w:rPr><w:t>a reference to bookmark1</w:t></w:r>' +
         '<w:r w:rsidR="' + $newRunRsid + '"><w:fldChar w:fldCharType="end"/></w:r>' +
         '<w:r w:rsidR="00D0546C"><w:t xml:space="preserve"> </w:t></w:r>' +
         '</w:p>'
$null = $p5.Range.InsertXML($p5Xml)

Write-Host "duplicatedBookmarks ids refreshed (2.0.2 -> 2.0.3)"
